$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "OA" drug to "MMP62" and fix its Step_dilution value
$ws.Range("A2").Value = "MMP62"
$ws.Range("D2").Value = 5

# Add a new row for the "MMAE" drug
$ws.Range("A4").Value = "MMAE"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 20000
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 1

# Update the selection to reflect where the cursor ended up after editing
$ws.Range("F5").Select()
